# Auto-generated Excel COM-interop script to update cryptos.xlsx data
# Applies per the commit: "Updated cryptos list on Thu Sep  5 14:52:04 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the data range so that numeric-looking strings
# (e.g. "1.00", "504.74") are preserved as text instead of being parsed into numbers,
# matching the original inline-string cell type. We restore the style afterwards so
# no residual style/format changes remain on the cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '56.736.86'
$ws.Range("E2").Value = '  -0.39%  '

$ws.Range("D3").Value = '2.380.59'
$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = '504.74'
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").Value = '132.71'
$ws.Range("E6").Value = '  +0.28%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("D8").Value = '0.552'
$ws.Range("E8").Value = '  -0.83%  '

$ws.Range("D9").Value = '2.386.50'
$ws.Range("E9").Value = '  -1.70%  '

$ws.Range("D10").Value = '0.0982'
$ws.Range("E10").Value = '  +1.09%  '

$ws.Range("E11").Value = '  +0.37%  '

$ws.Range("D12").Value = '0.332'
$ws.Range("E12").Value = '  +3.21%  '

$ws.Range("D13").Value = '4.65'
$ws.Range("E13").Value = '  -0.32%  '

$ws.Range("D14").Value = '2.814.82'
$ws.Range("E14").Value = '  -0.10%  '

$ws.Range("D15").Value = '56.707.01'
$ws.Range("E15").Value = '  -0.05%  '

$ws.Range("D16").Value = '21.66'
$ws.Range("E16").Value = '  -0.58%  '

$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").Value = '  -0.23%  '

$ws.Range("D18").Value = '2.399.38'
$ws.Range("E18").Value = '  +1.29%  '

$ws.Range("D19").Value = '10.08'
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").Value = '309.65'
$ws.Range("E20").Value = '  -0.97%  '

$ws.Range("D21").Value = '4.03'
$ws.Range("E21").Value = '  -0.93%  '

$ws.Range("D22").Value = '6.20'
$ws.Range("E22").Value = '  -3.26%  '

$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").Value = '65.18'
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("E25").Value = '  +0.58%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '0.151'
$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = '0.370'
$ws.Range("E27").Value = '  -2.54%  '

$ws.Range("D28").Value = '7.33'
$ws.Range("E28").Value = '  -1.54%  '

$ws.Range("D29").Value = '171.83'
$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("D30").Value = '0.0₃0725'
$ws.Range("E30").Value = '  -1.54%  '

$ws.Range("D31").Value = '1.65'
$ws.Range("E31").Value = '  -2.23%  '

$ws.Range("E32").Value = '  -2.83%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.21%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '5.80'
$ws.Range("E34").Value = '  -5.99%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").Value = '17.78'
$ws.Range("E36").Value = '  -1.30%  '

$ws.Range("D37").Value = '1.19'
$ws.Range("E37").Value = '  -2.77%  '

$ws.Range("D38").Value = '3.83'
$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("D39").Value = '0.816'
$ws.Range("E39").Value = '  +0.39%  '

$ws.Range("D40").Value = '36.15'
$ws.Range("E40").Value = '  +0.60%  '

$ws.Range("D41").Value = '1.43'
$ws.Range("E41").Value = '  -1.82%  '

$ws.Range("D42").Value = '131.03'
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("D43").Value = '3.38'
$ws.Range("E43").Value = '  -0.84%  '

$ws.Range("D44").Value = '4.82'
$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("D45").Value = '0.564'
$ws.Range("E45").Value = '  -0.39%  '

$ws.Range("D46").Value = '0.0910'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").Value = '247.58'
$ws.Range("E47").Value = '  -4.32%  '

$ws.Range("D48").Value = '0.0484'
$ws.Range("E48").Value = '  -2.09%  '

$ws.Range("D49").Value = '0.0209'
$ws.Range("E49").Value = '  -1.30%  '

$ws.Range("D50").Value = '17.14'
$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("D51").Value = '1.57'
$ws.Range("E51").Value = '  -1.20%  '

# Restore default (Normal) style on the data range so cells have no explicit
# number-format override left behind (matches original workbook formatting).
$dataRange.Style = "Normal"
